# Updates the Price (D) and Volume(1h) (E) columns of the cryptos sheet
# with freshly scraped values, per GitHub Actions automation run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new text value. Values that parse as plain numbers are prefixed
# with a literal apostrophe so Excel keeps them as text (matching the
# original formatting, e.g. '1.00' must not collapse to '1').
$updates = [ordered]@{
    # Row 2: Bitcoin
    D2 = '69.488.97'
    E2 = '  +1.65%  '
    # Row 3: Ethereum
    D3 = '3.944.62'
    E3 = '  +0.28%  '
    # Row 4: TetherUSD
    D4 = '''0.999'
    E4 = '  -0.10%  '
    # Row 5: BNB
    D5 = '''503.82'
    E5 = '  +3.35%  '
    # Row 6: Solana
    D6 = '''147.70'
    E6 = '  -0.61%  '
    # Row 7: XRP
    E7 = '  -0.53%  '
    # Row 8: USDC
    E8 = '  +0.00%  '
    # Row 9: Cardano
    D9 = '''0.734'
    E9 = '  -0.37%  '
    # Row 10: Dogecoin
    E10 = '  +3.97%  '
    # Row 11: ShibaInu
    D11 = '''0.0000350'
    # Row 12: Avalanche
    D12 = '''43.57'
    E12 = '  +1.00%  '
    # Row 13: Polkadot
    D13 = '''10.50'
    E13 = '  -2.05%  '
    # Row 14: WrappedliquidstakedEther2.0
    D14 = '4.574.86'
    E14 = '  +0.26%  '
    # Row 15: WrappedEther
    D15 = '3.950.97'
    E15 = '  +0.69%  '
    # Row 16: Uniswap
    D16 = '''14.24'
    E16 = '  -2.22%  '
    # Row 17: TRON
    E17 = '  -0.29%  '
    # Row 18: Polygon
    E18 = '  +5.03%  '
    # Row 19: Chainlink
    D19 = '''20.02'
    E19 = '  -0.09%  '
    # Row 20: WrappedBTC
    D20 = '69.485.67'
    E20 = '  +1.49%  '
    # Row 21: BitcoinCash
    D21 = '''436.66'
    E21 = '  -1.60%  '
    # Row 22: ImmutableX
    D22 = '''3.46'
    E22 = '  -2.04%  '
    # Row 23: InternetComputer(DFINITY)
    D23 = '''14.73'
    E23 = '  -2.80%  '
    # Row 24: Litecoin
    D24 = '''89.02'
    E24 = '  +0.49%  '
    # Row 25: RenderToken
    D25 = '''12.00'
    E25 = '  +4.11%  '
    # Row 26: PancakeSwap
    E26 = '  +6.45%  '
    # Row 27: Filecoin
    D27 = '''11.21'
    E27 = '  -2.16%  '
    # Row 28: EthereumClassic
    D28 = '''37.19'
    E28 = '  -4.68%  '
    # Row 29: LEO
    D29 = '''5.67'
    E29 = '  -3.04%  '
    # Row 30: Bittensor
    D30 = '''708.94'
    E30 = '  -1.82%  '
    # Row 31: Cosmos
    E31 = '  -1.91%  '
    # Row 32: Hedera
    E32 = '  -1.34%  '
    # Row 33: Toncoin
    E33 = '  -0.82%  '
    # Row 34: OKB
    D34 = '''64.59'
    E34 = '  +5.28%  '
    # Row 35: TheGraph
    E35 = '  +12.50%  '
    # Row 36: PEPE
    D36 = '0.0₃0896'
    E36 = '  -1.64%  '
    # Row 37: NEARProtocol
    D37 = '''6.05'
    E37 = '  -3.77%  '
    # Row 38: InjectiveProtocol
    D38 = '''40.99'
    E38 = '  -3.53%  '
    # Row 39: Kaspa
    E39 = '  +1.06%  '
    # Row 40: Dai
    D40 = '''1.00'
    E40 = '  -0.06%  '
    # Row 41: FirstDigitalUSD
    E41 = '  -0.02%  '
    # Row 42: VeChain
    D42 = '''0.0491'
    E42 = '  +1.63%  '
    # Row 43: Fetch.AI
    D43 = '''2.88'
    E43 = '  -5.80%  '
    # Row 44: ThetaToken
    D44 = '''3.08'
    E44 = '  -5.41%  '
    # Row 45: WEMIXToken
    D45 = '''3.07'
    E45 = '  +4.27%  '
    # Row 46: Stellar
    E46 = '  +0.92%  '
    # Row 47: ApeXProtocol
    D47 = '''3.38'
    E47 = '  +3.17%  '
    # Row 48: Stacks
    D48 = '''3.02'
    E48 = '  +5.86%  '
    # Row 49: LidoDAOToken
    D49 = '''3.41'
    E49 = '  -0.73%  '
    # Row 50: BabyDogeCoin
    D50 = '0.0₆0348'
    E50 = '  -5.36%  '
    # Row 51: ARBITRUM
    D51 = '''2.11'
    E51 = '  -2.66%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

